$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.556.81"
$ws.Range("E2").Value = "  +2.46%  "
$ws.Range("D3").Value = "1.642.45"
$ws.Range("E3").Value = "  +4.33%  "
$style_D4 = $ws.Range("D4").Style
$ws.Range("D4").Value = "'0.9980"
$ws.Range("D4").Style = $style_D4
$ws.Range("E4").Value = "  -0.37%  "
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").Value = "'308.20"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = "  +3.13%  "
$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.9978"
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = "  -0.49%  "
$style_D7 = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.3783"
$ws.Range("D7").Style = $style_D7
$ws.Range("E7").Value = "  +1.29%  "
$style_D8 = $ws.Range("D8").Style
$ws.Range("D8").Value = "'52.99"
$ws.Range("D8").Style = $style_D8
$ws.Range("E8").Value = "  +6.22%  "
$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.3679"
$ws.Range("D9").Style = $style_D9
$ws.Range("E9").Value = "  +3.86%  "
$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").Value = "'1.287"
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = "  +6.58%  "
$style_D11 = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.08200"
$ws.Range("D11").Style = $style_D11
$ws.Range("E11").Value = "  +3.48%  "
$style_D12 = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.9981"
$ws.Range("D12").Style = $style_D12
$ws.Range("E12").Value = "  -0.37%  "
$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").Value = "'23.32"
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = "  +7.38%  "
$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").Value = "'6.683"
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = "  +4.63%  "
$ws.Range("E15").Value = "  +6.25%  "
$style_D16 = $ws.Range("D16").Style
$ws.Range("D16").Value = "'7.490"
$ws.Range("D16").Style = $style_D16
$ws.Range("E16").Value = "  +3.57%  "
$ws.Range("D17").Value = "1.640.17"
$ws.Range("E17").Value = "  +3.95%  "
$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").Value = "'95.04"
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = "  +3.90%  "
$ws.Range("E19").Value = "  +3.26%  "
$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").Value = "'18.48"
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = "  +4.89%  "
$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").Value = "'6.609"
$ws.Range("D21").Style = $style_D21
$ws.Range("E21").Value = "  +4.38%  "
$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").Value = "'0.9975"
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").Value = "23.551.71"
$ws.Range("E23").Value = "  +2.45%  "
$style_D24 = $ws.Range("D24").Style
$ws.Range("D24").Value = "'13.02"
$ws.Range("D24").Style = $style_D24
$ws.Range("E24").Value = "  +3.69%  "
$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").Value = "'3.135"
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = "  +11.62%  "
$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").Value = "'2.428"
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = "  +2.95%  "
$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").Value = "'21.46"
$ws.Range("D27").Style = $style_D27
$ws.Range("E27").Value = "  +4.80%  "
$style_D28 = $ws.Range("D28").Style
$ws.Range("D28").Value = "'151.84"
$ws.Range("D28").Style = $style_D28
$ws.Range("E28").Value = "  +3.39%  "
$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").Value = "'5.320"
$ws.Range("D29").Style = $style_D29
$ws.Range("E29").Value = "  +3.17%  "
$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").Value = "'136.59"
$ws.Range("D30").Style = $style_D30
$ws.Range("E30").Value = "  +4.32%  "
$style_D31 = $ws.Range("D31").Style
$ws.Range("D31").Value = "'2.428"
$ws.Range("D31").Style = $style_D31
$ws.Range("E31").Value = "  +3.89%  "
$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'6.867"
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("D33").Value = "1.812.91"
$ws.Range("E33").Value = "  +3.54%  "
$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.9784"
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = "  +5.79%  "
$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").Value = "'0.02835"
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = "  +7.25%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'10.47"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = "  +6.13%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$style_D37 = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.07505"
$ws.Range("D37").Style = $style_D37
$ws.Range("E37").Value = "  +3.01%  "
$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").Value = "'6.243"
$ws.Range("D38").Style = $style_D38
$ws.Range("E38").Value = "  +5.05%  "
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.2547"
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = "  +3.89%  "
$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.08871"
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = "  +1.81%  "
$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'1.402"
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = "  +5.03%  "
$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.7186"
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = "  +5.34%  "
$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").Value = "'12.73"
$ws.Range("D43").Style = $style_D43
$ws.Range("E43").Value = "  +8.45%  "
$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").Value = "'16.20"
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = "  +11.61%  "
$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.6658"
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = "  +5.91%  "
$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").Value = "'2.372"
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = "  +6.34%  "
$ws.Range("E47").Value = "  +2.18%  "
$style_D48 = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.9963"
$ws.Range("D48").Style = $style_D48
$ws.Range("E48").Value = "  -0.52%  "
$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.08070"
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = "  +2.97%  "
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").Value = "'131.78"
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("E51").Value = "  +3.01%  "
